$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Thbs1"
$ws.Range("C2").Value = "Tnfrsf11b"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 149.656361
$ws.Range("H2").Value = 448.969083
$ws.Range("I2").Value = 0.5921360794347563
$ws.Range("J2").Value = 0.5921360794347564
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.095195666666667
$ws.Range("N2").Value = 6.285587
$ws.Range("O2").Value = 0.8546922300706357
$ws.Range("P2").Value = 0.8546922300706358
$ws.Range("Q2").Value = 313.5593590563023
$ws.Range("R2").Value = 2822.034231506721
$ws.Range("S2").Value = 0.5060941062373749
$ws.Range("T2").Value = 0.5060941062373751

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Thbs1"
$ws.Range("C3").Value = "Tnfrsf11b"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 149.656361
$ws.Range("H3").Value = 448.969083
$ws.Range("I3").Value = 0.5921360794347563
$ws.Range("J3").Value = 0.5921360794347564
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.356208
$ws.Range("N3").Value = 1.068624
$ws.Range("O3").Value = 0.1453077699293643
$ws.Range("P3").Value = 0.1453077699293643
$ws.Range("Q3").Value = 53.30879303908799
$ws.Range("R3").Value = 479.779137351792
$ws.Range("S3").Value = 0.08604197319738134
$ws.Range("T3").Value = 0.08604197319738136

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Thbs1"
$ws.Range("C4").Value = "Tnfrsf11b"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 52.73412466666667
$ws.Range("H4").Value = 158.202374
$ws.Range("I4").Value = 0.208649853730866
$ws.Range("J4").Value = 0.208649853730866
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.095195666666667
$ws.Range("N4").Value = 6.285587
$ws.Range("O4").Value = 0.8546922300706357
$ws.Range("P4").Value = 0.8546922300706358
$ws.Range("Q4").Value = 110.4883094870598
$ws.Range("R4").Value = 994.3947853835379
$ws.Range("S4").Value = 0.1783314087891458
$ws.Range("T4").Value = 0.1783314087891459

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Thbs1"
$ws.Range("C5").Value = "Tnfrsf11b"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 52.73412466666667
$ws.Range("H5").Value = 158.202374
$ws.Range("I5").Value = 0.208649853730866
$ws.Range("J5").Value = 0.208649853730866
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.356208
$ws.Range("N5").Value = 1.068624
$ws.Range("O5").Value = 0.1453077699293643
$ws.Range("P5").Value = 0.1453077699293643
$ws.Range("Q5").Value = 18.784317079264
$ws.Range("R5").Value = 169.058853713376
$ws.Range("S5").Value = 0.03031844494172019
$ws.Range("T5").Value = 0.03031844494172019

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Thbs1"
$ws.Range("C6").Value = "Tnfrsf11b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 50.34932566666667
$ws.Range("H6").Value = 151.047977
$ws.Range("I6").Value = 0.1992140668343777
$ws.Range("J6").Value = 0.1992140668343777
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.095195666666667
$ws.Range("N6").Value = 6.285587
$ws.Range("O6").Value = 0.8546922300706357
$ws.Range("P6").Value = 0.8546922300706358
$ws.Range("Q6").Value = 105.4916889563888
$ws.Range("R6").Value = 949.4252006074989
$ws.Range("S6").Value = 0.1702667150441149
$ws.Range("T6").Value = 0.1702667150441149

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Thbs1"
$ws.Range("C7").Value = "Tnfrsf11b"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 50.34932566666667
$ws.Range("H7").Value = 151.047977
$ws.Range("I7").Value = 0.1992140668343777
$ws.Range("J7").Value = 0.1992140668343777
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.356208
$ws.Range("N7").Value = 1.068624
$ws.Range("O7").Value = 0.1453077699293643
$ws.Range("P7").Value = 0.1453077699293643
$ws.Range("Q7").Value = 17.934832597072
$ws.Range("R7").Value = 161.413493373648
$ws.Range("S7").Value = 0.02894735179026275
$ws.Range("T7").Value = 0.02894735179026275
